$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "target" column header (O1), reuses existing header cell style
$ws.Range("O1").Value = "target"

# Update data grid (rows 2-18, columns A-O) with spectral-clustering results
$ws.Cells.Item(2,1).Value = 10
$ws.Cells.Item(2,2).Value = 1.068282978454066
$ws.Cells.Item(2,3).Value = 1.086827571115874
$ws.Cells.Item(2,4).Value = -0.2830771045654332
$ws.Cells.Item(2,5).Value = -0.5775421200673674
$ws.Cells.Item(2,6).Value = 0.1077360276293435
$ws.Cells.Item(2,7).Value = 0.1041029176906046
$ws.Cells.Item(2,8).Value = 0.1923087798717137
$ws.Cells.Item(2,9).Value = 0.2894757767297265
$ws.Cells.Item(2,10).Value = 0.5166036693310342
$ws.Cells.Item(2,11).Value = 0.5516134326604308
$ws.Cells.Item(2,12).Value = -0.5481823527185627
$ws.Cells.Item(2,13).Value = -0.4214390128295122
$ws.Cells.Item(2,14).Value = -0.2891866838178406
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(3,1).Value = 32
$ws.Cells.Item(3,2).Value = -0.7170444859158436
$ws.Cells.Item(3,3).Value = -0.9626627575504668
$ws.Cells.Item(3,4).Value = -1.326370565969052
$ws.Cells.Item(3,5).Value = -0.3508826171547983
$ws.Cells.Item(3,6).Value = -0.1116572877339215
$ws.Cells.Item(3,7).Value = 0.3025477387642772
$ws.Cells.Item(3,8).Value = 0.02227960656881003
$ws.Cells.Item(3,9).Value = -0.8292441429574355
$ws.Cells.Item(3,10).Value = -0.3231513895542691
$ws.Cells.Item(3,11).Value = -0.6051667266971232
$ws.Cells.Item(3,12).Value = 0.4240250088644789
$ws.Cells.Item(3,13).Value = 0.5689230922157072
$ws.Cells.Item(3,14).Value = 0.5941888100937924
$ws.Cells.Item(3,15).Value = 1
$ws.Cells.Item(4,1).Value = 1
$ws.Cells.Item(4,2).Value = 1.083936345210007
$ws.Cells.Item(4,3).Value = 1.093264700238166
$ws.Cells.Item(4,4).Value = 0.1483705827004671
$ws.Cells.Item(4,5).Value = -0.2229674175134839
$ws.Cells.Item(4,6).Value = -0.4549054893864045
$ws.Cells.Item(4,7).Value = 0.213315262424153
$ws.Cells.Item(4,8).Value = -0.5341985581575719
$ws.Cells.Item(4,9).Value = -0.0575880076469479
$ws.Cells.Item(4,10).Value = -0.3098117486188885
$ws.Cells.Item(4,11).Value = -0.9285690794037499
$ws.Cells.Item(4,12).Value = -1.478019162058223
$ws.Cells.Item(4,13).Value = -1.50507127255204
$ws.Cells.Item(4,14).Value = -1.83833156063562
$ws.Cells.Item(4,15).Value = 1
$ws.Cells.Item(5,1).Value = 9
$ws.Cells.Item(5,2).Value = -0.9950186759999999
$ws.Cells.Item(5,3).Value = -1.14275891
$ws.Cells.Item(5,4).Value = -0.8379565390000001
$ws.Cells.Item(5,5).Value = -1.356096241
$ws.Cells.Item(5,6).Value = -1.482855771
$ws.Cells.Item(5,7).Value = -0.116152378
$ws.Cells.Item(5,8).Value = -1.960741118
$ws.Cells.Item(5,9).Value = -2.808609729
$ws.Cells.Item(5,10).Value = -1.937638123
$ws.Cells.Item(5,11).Value = -2.478579661
$ws.Cells.Item(5,12).Value = -1.538677731
$ws.Cells.Item(5,13).Value = -0.954723074
$ws.Cells.Item(5,14).Value = -0.998195754
$ws.Cells.Item(5,15).Value = 1
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = 0.078267607
$ws.Cells.Item(6,3).Value = 0.147016056
$ws.Cells.Item(6,4).Value = -0.099193434
$ws.Cells.Item(6,5).Value = -0.9932849579999999
$ws.Cells.Item(6,6).Value = -0.36166605
$ws.Cells.Item(6,7).Value = 0.204833342
$ws.Cells.Item(6,8).Value = -0.837582409
$ws.Cells.Item(6,9).Value = -0.6419876839999999
$ws.Cells.Item(6,10).Value = -0.254815942
$ws.Cells.Item(6,11).Value = -0.379259574
$ws.Cells.Item(6,12).Value = -1.062085395
$ws.Cells.Item(6,13).Value = -1.581943415
$ws.Cells.Item(6,14).Value = -1.785023628
$ws.Cells.Item(6,15).Value = 1
$ws.Cells.Item(7,1).Value = 8
$ws.Cells.Item(7,2).Value = -1.113243525
$ws.Cells.Item(7,3).Value = -1.218439747
$ws.Cells.Item(7,4).Value = -0.46169408
$ws.Cells.Item(7,5).Value = -0.506681916
$ws.Cells.Item(7,6).Value = -0.141881335
$ws.Cells.Item(7,7).Value = 0.22744881
$ws.Cells.Item(7,8).Value = -0.997138382
$ws.Cells.Item(7,9).Value = -1.313530996
$ws.Cells.Item(7,10).Value = -0.829717461
$ws.Cells.Item(7,11).Value = -0.556270676
$ws.Cells.Item(7,12).Value = -0.016803582
$ws.Cells.Item(7,13).Value = -0.136466835
$ws.Cells.Item(7,14).Value = 0.219988164
$ws.Cells.Item(7,15).Value = 1
$ws.Cells.Item(8,1).Value = 29
$ws.Cells.Item(8,2).Value = -1.270561510282354
$ws.Cells.Item(8,3).Value = -1.123748824026904
$ws.Cells.Item(8,4).Value = -0.09099933402576334
$ws.Cells.Item(8,5).Value = -0.04498554308108361
$ws.Cells.Item(8,6).Value = 0.1352646616835821
$ws.Cells.Item(8,7).Value = 0.6061199888280842
$ws.Cells.Item(8,8).Value = 0.5865008059729586
$ws.Cells.Item(8,9).Value = 0.1804280281973176
$ws.Cells.Item(8,10).Value = -0.6615487182621481
$ws.Cells.Item(8,11).Value = 0.4701499509641796
$ws.Cells.Item(8,12).Value = 1.121854132091655
$ws.Cells.Item(8,13).Value = 1.101580566908483
$ws.Cells.Item(8,14).Value = 1.146047200980204
$ws.Cells.Item(8,15).Value = 1
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = -1.162809267
$ws.Cells.Item(9,3).Value = -1.230742598
$ws.Cells.Item(9,4).Value = -0.481548895
$ws.Cells.Item(9,5).Value = -0.48584958
$ws.Cells.Item(9,6).Value = -0.43719234
$ws.Cells.Item(9,7).Value = 0.114115978
$ws.Cells.Item(9,8).Value = -0.77417708
$ws.Cells.Item(9,9).Value = -1.367703306
$ws.Cells.Item(9,10).Value = -0.7770090820000001
$ws.Cells.Item(9,11).Value = -0.805050077
$ws.Cells.Item(9,12).Value = -0.044871051
$ws.Cells.Item(9,13).Value = 0.152162831
$ws.Cells.Item(9,14).Value = 0.264317159
$ws.Cells.Item(9,15).Value = 1
$ws.Cells.Item(10,1).Value = 50
$ws.Cells.Item(10,2).Value = -0.8758349393820211
$ws.Cells.Item(10,3).Value = -0.7819755033952668
$ws.Cells.Item(10,4).Value = -0.1136147965221389
$ws.Cells.Item(10,5).Value = -0.185492660811697
$ws.Cells.Item(10,6).Value = -0.0781696765984117
$ws.Cells.Item(10,7).Value = -0.7311395274105609
$ws.Cells.Item(10,8).Value = -0.7733672661317547
$ws.Cells.Item(10,9).Value = 0.2632672729876283
$ws.Cells.Item(10,10).Value = -0.09239871205352695
$ws.Cells.Item(10,11).Value = 0.04069870965258254
$ws.Cells.Item(10,12).Value = -0.3480279619857963
$ws.Cells.Item(10,13).Value = -0.1790791255580886
$ws.Cells.Item(10,14).Value = 0.1553656948685382
$ws.Cells.Item(10,15).Value = 2
$ws.Cells.Item(11,1).Value = 46
$ws.Cells.Item(11,2).Value = -0.554723695977998
$ws.Cells.Item(11,3).Value = -0.7348956677261104
$ws.Cells.Item(11,4).Value = -2.200403713486279
$ws.Cells.Item(11,5).Value = -0.7021511342036996
$ws.Cells.Item(11,6).Value = -1.10566236468863
$ws.Cells.Item(11,7).Value = -1.308502522194651
$ws.Cells.Item(11,8).Value = -0.4420220964019494
$ws.Cells.Item(11,9).Value = -0.8910558932634112
$ws.Cells.Item(11,10).Value = -0.3157577963177772
$ws.Cells.Item(11,11).Value = 0.03474938647192169
$ws.Cells.Item(11,12).Value = -0.03916370126099611
$ws.Cells.Item(11,13).Value = 0.03075681492401491
$ws.Cells.Item(11,14).Value = 0.1129110901943352
$ws.Cells.Item(11,15).Value = 2
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 0.08650412199999999
$ws.Cells.Item(12,3).Value = 0.16076415
$ws.Cells.Item(12,4).Value = 1.287529653
$ws.Cells.Item(12,5).Value = 1.458507674
$ws.Cells.Item(12,6).Value = 1.708264452
$ws.Cells.Item(12,7).Value = 1.737302749
$ws.Cells.Item(12,8).Value = 1.454654867
$ws.Cells.Item(12,9).Value = 1.685926169
$ws.Cells.Item(12,10).Value = 1.198936477
$ws.Cells.Item(12,11).Value = 1.009485194
$ws.Cells.Item(12,12).Value = 0.6228660539999999
$ws.Cells.Item(12,13).Value = 0.391809069
$ws.Cells.Item(12,14).Value = 0.180105296
$ws.Cells.Item(12,15).Value = 2
$ws.Cells.Item(13,1).Value = 14
$ws.Cells.Item(13,2).Value = -0.9819434020000001
$ws.Cells.Item(13,3).Value = -1.027168868
$ws.Cells.Item(13,4).Value = -0.391420759
$ws.Cells.Item(13,5).Value = -1.071957494
$ws.Cells.Item(13,6).Value = -1.015058936
$ws.Cells.Item(13,7).Value = -1.251488744
$ws.Cells.Item(13,8).Value = -1.889056123
$ws.Cells.Item(13,9).Value = -0.649593132
$ws.Cells.Item(13,10).Value = -0.321722455
$ws.Cells.Item(13,11).Value = -0.310259495
$ws.Cells.Item(13,12).Value = -0.362302266
$ws.Cells.Item(13,13).Value = -0.24811692
$ws.Cells.Item(13,14).Value = 0.142904433
$ws.Cells.Item(13,15).Value = 2
$ws.Cells.Item(14,1).Value = 43
$ws.Cells.Item(14,2).Value = -0.3543593665836199
$ws.Cells.Item(14,3).Value = -0.607770836055068
$ws.Cells.Item(14,4).Value = -1.023579362476468
$ws.Cells.Item(14,5).Value = -0.955669614345772
$ws.Cells.Item(14,6).Value = -1.000502683286009
$ws.Cells.Item(14,7).Value = -1.440570144316475
$ws.Cells.Item(14,8).Value = -0.6574608486657882
$ws.Cells.Item(14,9).Value = -0.3437044833041243
$ws.Cells.Item(14,10).Value = -0.3136140931433612
$ws.Cells.Item(14,11).Value = 0.1547769421110219
$ws.Cells.Item(14,12).Value = 0.2203860173098343
$ws.Cells.Item(14,13).Value = 0.1150587133144156
$ws.Cells.Item(14,14).Value = -0.0393656080207815
$ws.Cells.Item(14,15).Value = 2
$ws.Cells.Item(15,1).Value = 16
$ws.Cells.Item(15,2).Value = -0.898339942
$ws.Cells.Item(15,3).Value = -0.9166695420000001
$ws.Cells.Item(15,4).Value = -0.124501255
$ws.Cells.Item(15,5).Value = 0.389349156
$ws.Cells.Item(15,6).Value = 1.017939856
$ws.Cells.Item(15,7).Value = 1.065573911
$ws.Cells.Item(15,8).Value = 1.109928244
$ws.Cells.Item(15,9).Value = 0.875583775
$ws.Cells.Item(15,10).Value = 0.875780095
$ws.Cells.Item(15,11).Value = 1.330532029
$ws.Cells.Item(15,12).Value = 1.080667858
$ws.Cells.Item(15,13).Value = 0.87872014
$ws.Cells.Item(15,14).Value = 0.4082125
$ws.Cells.Item(15,15).Value = 3
$ws.Cells.Item(16,1).Value = 59
$ws.Cells.Item(16,2).Value = 1.098824281425039
$ws.Cells.Item(16,3).Value = 1.124544984094857
$ws.Cells.Item(16,4).Value = 1.730050599012491
$ws.Cells.Item(16,5).Value = 1.85323807668174
$ws.Cells.Item(16,6).Value = 2.19314676580624
$ws.Cells.Item(16,7).Value = 2.024663755796917
$ws.Cells.Item(16,8).Value = 2.554541131790494
$ws.Cells.Item(16,9).Value = 2.318109608543673
$ws.Cells.Item(16,10).Value = 2.100052046017395
$ws.Cells.Item(16,11).Value = 2.044557233732532
$ws.Cells.Item(16,12).Value = 1.594573606845447
$ws.Cells.Item(16,13).Value = 1.329427998478601
$ws.Cells.Item(16,14).Value = 1.092029141150236
$ws.Cells.Item(16,15).Value = 3
$ws.Cells.Item(17,1).Value = 62
$ws.Cells.Item(17,2).Value = -0.641595727756713
$ws.Cells.Item(17,3).Value = -0.5590221552278335
$ws.Cells.Item(17,4).Value = 0.3519448183223575
$ws.Cells.Item(17,5).Value = 1.25896503686056
$ws.Cells.Item(17,6).Value = 1.622979822911583
$ws.Cells.Item(17,7).Value = 0.5815436039527361
$ws.Cells.Item(17,8).Value = 1.597246989459927
$ws.Cells.Item(17,9).Value = 1.191527893722771
$ws.Cells.Item(17,10).Value = 1.227082251603201
$ws.Cells.Item(17,11).Value = 1.437210920502419
$ws.Cells.Item(17,12).Value = 1.106532648937356
$ws.Cells.Item(17,13).Value = 0.6944463772260019
$ws.Cells.Item(17,14).Value = 0.6524015959985536
$ws.Cells.Item(17,15).Value = 3
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = 0.080592808
$ws.Cells.Item(18,3).Value = 0.154079561
$ws.Cells.Item(18,4).Value = 0.950204087
$ws.Cells.Item(18,5).Value = 0.5333746429999999
$ws.Cells.Item(18,6).Value = 0.669714496
$ws.Cells.Item(18,7).Value = 1.32170149
$ws.Cells.Item(18,8).Value = 0.394812815
$ws.Cells.Item(18,9).Value = 0.060439413
$ws.Cells.Item(18,10).Value = 0.590968691
$ws.Cells.Item(18,11).Value = 0.584491591
$ws.Cells.Item(18,12).Value = 0.7292046520000001
$ws.Cells.Item(18,13).Value = 0.752660914
$ws.Cells.Item(18,14).Value = 0.4082125
$ws.Cells.Item(18,15).Value = 3
